$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 0.6219589999999999
$ws.Range("I2").Value = 0.6398583988494134
$ws.Range("J2").Value = 0.6398583988494134
$ws.Range("S2").Value = 0.6398583988494134
$ws.Range("T2").Value = 0.6398583988494134

# Row 3 updates
$ws.Range("G3").Value = 0.116689
$ws.Range("H3").Value = 0.350067
$ws.Range("I3").Value = 0.3601416011505865
$ws.Range("J3").Value = 0.3601416011505865
$ws.Range("Q3").Value = 0.06381589162466667
$ws.Range("R3").Value = 0.574343024622
$ws.Range("S3").Value = 0.3601416011505865
$ws.Range("T3").Value = 0.3601416011505865
